# SCI 1402 Time Keeping.xlsx - add two new logged tasks (rows 15-16) and
# extend the "Running Total (h)" formula (column D) down through row 32,
# matching the author's fill-down of the running-total column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: "Processing datasets" (6h on 3/2/2025) -----------------------
# Copy the date formatting (style) from A13 so the new date cells keep the
# existing short-date number format instead of minting a new one.
$ws.Range("A13").Copy()
$ws.Range("A15:A16").PasteSpecial(-4122)

$ws.Range("A15").Value = 45718
$ws.Range("B15").Value = "Processing datasets"
$ws.Range("C15").Value = 6
$ws.Range("D15").Formula = "=SUM(INDEX(C:C,1):INDEX(C:C,ROW()))"

# --- Row 16: "Merging and cleaning the combined dataset" (8h on 3/2/2025) -
$ws.Range("A16").Value = 45718
$ws.Range("B16").Value = "Merging and cleaning the combined dataset"
$ws.Range("C16").Value = 8
$ws.Range("D16").Formula = "=SUM(INDEX(C:C,1):INDEX(C:C,ROW()))"

# --- Rows 17-32: running-total formula filled down, no other data yet ----
for ($r = 17; $r -le 32; $r++) {
    $ws.Range("D$r").Formula = "=SUM(INDEX(C:C,1):INDEX(C:C,ROW()))"
}

# Selection left where the author's cursor ended up after filling down.
$ws.Range("D17").Select()
